# Apply the latest cryptos list update (GitHub Actions refresh).
# All cells in this sheet are plain text (inline/shared strings), including
# price columns like "35.136.16" or "42.01" that look numeric. Setting
# NumberFormat to text ("@") before the assignment stops Excel from silently
# re-typing them as numbers, and ClearFormats() afterwards drops that helper
# format again so the cell keeps its original (default) style - only the text
# content changes, exactly like the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "35.136.16"
Set-TextValue "E2" "  +0.90%  "
Set-TextValue "D3" "1.857.07"
Set-TextValue "E3" "  +1.77%  "
Set-TextValue "E4" "  +0.26%  "
Set-TextValue "D5" "238.97"
Set-TextValue "E5" "  +3.72%  "
Set-TextValue "E6" "  +1.21%  "
Set-TextValue "E7" "  +0.26%  "
Set-TextValue "D8" "42.01"
Set-TextValue "E8" "  +5.83%  "
Set-TextValue "D9" "0.329"
Set-TextValue "E9" "  +2.96%  "
Set-TextValue "E10" "  +1.64%  "
Set-TextValue "E11" "  +0.10%  "
Set-TextValue "D12" "2.124.60"
Set-TextValue "E12" "  +1.66%  "
Set-TextValue "D13" "11.48"
Set-TextValue "E13" "  +1.96%  "
Set-TextValue "B14" "Polygon"
Set-TextValue "C14" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.677"
Set-TextValue "E14" "  +1.73%  "
Set-TextValue "B15" "WrappedEther"
Set-TextValue "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D15" "1.834.28"
Set-TextValue "E15" "  +0.58%  "
Set-TextValue "E16" "  +1.81%  "
Set-TextValue "D17" "35.104.81"
Set-TextValue "E17" "  +1.11%  "
Set-TextValue "D18" "69.84"
Set-TextValue "E18" "  +0.63%  "
Set-TextValue "D19" "0.0₃0795"
Set-TextValue "E19" "  +1.33%  "
Set-TextValue "D20" "240.81"
Set-TextValue "E20" "  +0.60%  "
Set-TextValue "D21" "12.26"
Set-TextValue "E21" "  +1.41%  "
Set-TextValue "E22" "  +1.61%  "
Set-TextValue "E23" "  +0.36%  "
Set-TextValue "E24" "  +0.37%  "
Set-TextValue "D25" "167.76"
Set-TextValue "E25" "  -3.31%  "
Set-TextValue "D26" "1.89"
Set-TextValue "E26" "  +26.65%  "
Set-TextValue "D27" "7.97"
Set-TextValue "E27" "  +3.51%  "
Set-TextValue "E28" "  +2.16%  "
Set-TextValue "E29" "  +0.34%  "
Set-TextValue "E30" "  +0.26%  "
Set-TextValue "E31" "  +1.48%  "
Set-TextValue "E32" "  +2.28%  "
Set-TextValue "E33" "  +27.42%  "
Set-TextValue "E34" "  +2.63%  "
Set-TextValue "D35" "0.835"
Set-TextValue "E35" "  +19.90%  "
Set-TextValue "E36" "  +11.85%  "
Set-TextValue "E37" "  +7.24%  "
Set-TextValue "E38" "  +7.71%  "
Set-TextValue "D39" "90.59"
Set-TextValue "E39" "  -1.30%  "
Set-TextValue "E40" "  +4.26%  "
Set-TextValue "D41" "1.342.13"
Set-TextValue "E41" "  +0.37%  "
Set-TextValue "D42" "14.93"
Set-TextValue "E42" "  +3.60%  "
Set-TextValue "E43" "  +3.76%  "
Set-TextValue "B44" "HuobiToken"
Set-TextValue "C44" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D44" "2.41"
Set-TextValue "E44" "  -0.54%  "
Set-TextValue "B45" "Gas"
Set-TextValue "C45" "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
Set-TextValue "D45" "12.49"
Set-TextValue "E45" "  +44.78%  "
Set-TextValue "D46" "0.0556"
Set-TextValue "E46" "  +6.49%  "
Set-TextValue "E47" "  -0.26%  "
Set-TextValue "D48" "6.62"
Set-TextValue "E48" "  +5.79%  "
Set-TextValue "D49" "2.040.29"
Set-TextValue "E49" "  +1.61%  "
Set-TextValue "E50" "  +1.58%  "
Set-TextValue "E51" "  +0.34%  "
